# POCOR-4156 - report_card_template.xlsx
# Rename the body-mass report-card placeholders from the old
# "body_mass" association to the new "student_body_mass" one
# (StudentBodyMassesController / UserBodyMassesTable rename).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("General")

$ws.Range("E9").Value = '${InstitutionStudentsReportCards.student_body_mass.height} m'
$ws.Range("E10").Value = '${InstitutionStudentsReportCards.student_body_mass.weight} kg'
$ws.Range("E11").Value = '${InstitutionStudentsReportCards.student_body_mass.body_mass_index}'

$ws.Range("E11").Select()
